$d = $word.ActiveDocument

# 1) Ativacao date bump
$d.Content.Find.Execute("Ativação: 01/01/2017", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2) | Out-Null

# 2) Insert two new docente lines before Daniela (as separate runs with their own line breaks)
$r2 = $d.Content
$r2.Find.Execute("5840963 - Daniela Camargo Vernilli", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r2.InsertBefore("3586455 - Cassius Olivio Figueiredo Terra Ruchert`v")
$r2.InsertBefore("144651 - Antonio Fernando Sartori`v")

# 3) Replace Gustavo line with Hugo line
$d.Content.Find.Execute("5840820 - Gustavo Aristides Santana Martinez", $true, $false, $false, $false, $false, $true, 1, $false, "984972 - Hugo Ricardo Zschommler Sandim", 2) | Out-Null

# 4) Programa resumido paragraph
$d.Content.Find.Execute("Introdução aos materiais para a indústria química e outras, propriedades,  especificações, seleção, fabricação, aplicação e possíveis falhas.", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução aos materiais para a indústria química, propriedades, especificações, seleção, fabricação, aplicação e corrosão.", 2) | Out-Null

# 5) Programa (long) paragraph
$d.Content.Find.Execute("Introdução aos materiais. - Seleção de materiais. - Fatores que influenciam na seleção dos materiais (indústria química, petroquímica, Nuclear e outras), melhoria das propriedades mecânicas dos metais. - Falhas em serviço e em processo.  Produtos siderúrgicos para aplicação em indústrias químicas - Aços carbono e especiais - Ferro fundido. – Processo de fabricação de aços e ferros fundidos, especificações, propriedades e aplicações.  Metais e ligas não ferrosas: especificações, propriedades e aplicações. Introdução à corrosão. - Causas e formas de corrosão. Proteção de superfícies metálicas contra a corrosão, revestimentos. Requisitos específicos de materiais metálicos para a indústria de óleo e gás.  Materiais não metálicos. Especificações, propriedades e aplicações.", $true, $false, $false, $false, $false, $true, 1, $false, "Introdução aos materiais. - Seleção de materiais. - Fatores que influenciam na seleção dos materiais (indústria química, petroquímica, Nuclear e outras), melhoria das propriedades mecânicas dos materiais. - Falhas em serviço e em processo.  Produtos siderúrgicos para aplicação em indústrias químicas - Aços carbono e especiais - Ferro fundido. – Processo de fabricação de aços e ferros fundidos, especificações, propriedades e aplicações. Metais e ligas não ferrosas e não metálicas: especificações, propriedades e aplicações. Introdução à corrosão aplicada a engenharia. Pilha Eletroquímica e eletrolítica, meios corrosivos, causas e formas de corrosão, corrosão seletiva, induzida por micromecanismos (MIQ), puntiforme, filiforme, frestas, CST etc... Proteção de superfícies metálicas contra a corrosão, tipo de revestimentos como aspersão térmica, PVD, QVD, etc..", 2) | Out-Null

# 6) Avaliacao - Metodo
$d.Content.Find.Execute("Duas provas", $true, $false, $false, $false, $false, $true, 1, $false, "De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.", 2) | Out-Null

# 7) Avaliacao - Criterio
$d.Content.Find.Execute("Serão aplicadas duas avaliações (P1 e P2) que comporão a nota final (NF). A nota final será calculada através da expressão: NF= (P1+P2)/2", $true, $false, $false, $false, $false, $true, 1, $false, "A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)", 2) | Out-Null

# 8) Avaliacao - Norma de recuperacao
$d.Content.Find.Execute("Para o aluno que obtiver Nota Final maior ou igual a 3,0 e menor do que 5,0, será aplicada uma avaliação de recuperação (NR), com pontuação de 0 a 10, que levará ao cálculo da média final(MF) através da seguinte expressão: MF=(NF+NR)/2.onde: NF=Nota Final e NR=Nota da Prova de Recuperação.", $true, $false, $false, $false, $false, $true, 1, $false, "Devido a cunho prático da disciplina não haverá recuperação.", 2) | Out-Null

